$d = $word.ActiveDocument

# 1) Delete the now-empty footnote (its body text is just a lone shad "།").
#    Footnote.Delete() removes both the <w:footnote> definition and the
#    <w:footnoteReference> run left in the body paragraph.
for ($i = $d.Footnotes.Count; $i -ge 1; $i--) {
    $fn = $d.Footnotes.Item($i)
    if ($fn.Range.Text -eq "།") {
        $fn.Delete()
    }
}

# 2) Append the folio marker to the text run that used to sit right before
#    that footnote reference.
$d.Content.Find.Execute("མཛད་པའོ།།", $true, $false, $false, $false, $false, $true, 1, $false, "མཛད་པའོ།།[༢༦༩བ]", 2)

# 3) Strip the stray "aa" typo left at the end of footnote 26's body text.
for ($i = 1; $i -le $d.Footnotes.Count; $i++) {
    $fn = $d.Footnotes.Item($i)
    $t = $fn.Range.Text
    if ($t.EndsWith("aa")) {
        $fn.Range.Text = $t.Substring(0, $t.Length - 2)
    }
}
